$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.13570000000001
$ws.Range("D4").Value = -8.189300000000005
$ws.Range("D5").Value = -8.761699999999996
$ws.Range("A6").Value = -20.18179999999999
$ws.Range("D6").Value = -8.498399999999997
$ws.Range("A7").Value = -21.29790000000001
$ws.Range("A8").Value = -20.65229999999999
$ws.Range("D8").Value = -8.346500000000002
$ws.Range("A16").Value = -20.30719999999999
$ws.Range("D16").Value = -8.248400000000004
$ws.Range("A20").Value = -22.19360000000002
$ws.Range("A21").Value = -20.49449999999998
$ws.Range("D22").Value = -8.153599999999997
